$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1800.75
$ws.Range("I33").Value = 1850.5
$ws.Range("J33").Value = 1751
$ws.Range("K33").Value = 1850.5
$ws.Range("L33").Value = 1751
$ws.Range("M33").Value = -1621.5
$ws.Range("N33").Value = -2209
$ws.Range("H80").Value = 38461750
$ws.Range("I80").Value = 50000120
$ws.Range("J80").Value = 532.3333
$ws.Range("K80").Value = 150000360
$ws.Range("L80").Value = 1596.9999
$ws.Range("M80").Value = -149999362
$ws.Range("N80").Value = -3592.9999
$ws.Range("H83").Value = 38461750
$ws.Range("I83").Value = 50000120
$ws.Range("J83").Value = 532.3333
$ws.Range("K83").Value = 450001080
$ws.Range("L83").Value = 4790.9997
$ws.Range("M83").Value = -449996088
$ws.Range("N83").Value = -14774.9997
$ws.Range("H94").Value = 14149.75
$ws.Range("I94").Value = 14149.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 14149.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -13698.75
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 6464
$ws.Range("I113").Value = 5400.5713
$ws.Range("J113").Value = 8325
$ws.Range("K113").Value = 5400.5713
$ws.Range("L113").Value = 8325
$ws.Range("M113").Value = -2146.5713
$ws.Range("N113").Value = -14833
$ws.Range("H127").Value = 4779.25
$ws.Range("I127").Value = 2900
$ws.Range("J127").Value = 5405.6665
$ws.Range("K127").Value = 8700
$ws.Range("L127").Value = 16216.9995
$ws.Range("M127").Value = -3740
$ws.Range("N127").Value = -26136.9995
$ws.Range("H138").Value = 4026.1052
$ws.Range("J138").Value = 4661.5386
$ws.Range("L138").Value = 13984.6158
$ws.Range("N138").Value = -24264.6158
$ws.Range("H141").Value = 4007.0476
$ws.Range("I141").Value = 4007.0476
$ws.Range("K141").Value = 12021.1428
$ws.Range("M141").Value = -6841.1428

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2087.2222
$ws.Range("I2").Value = 1338.0834
$ws.Range("K2").Value = 1338.0834
$ws.Range("M2").Value = -1225.0834
$ws.Range("H16").Value = 16413
$ws.Range("I16").Value = 33568.332
$ws.Range("J16").Value = 6119.8
$ws.Range("K16").Value = 33568.332
$ws.Range("L16").Value = 6119.8
$ws.Range("M16").Value = -33281.332
$ws.Range("N16").Value = -6693.8
$ws.Range("H32").Value = 23081
$ws.Range("I32").Value = 13834.479
$ws.Range("K32").Value = 13834.479
$ws.Range("M32").Value = -13547.479
$ws.Range("H74").Value = 2101.0476
$ws.Range("I74").Value = 2006.9474
$ws.Range("K74").Value = 2006.9474
$ws.Range("M74").Value = -1132.9474
$ws.Range("H77").Value = 2101.0476
$ws.Range("I77").Value = 2006.9474
$ws.Range("K77").Value = 10034.737
$ws.Range("M77").Value = -5666.737000000001
$ws.Range("H97").Value = 1404.1428
$ws.Range("I97").Value = 1404.1428
$ws.Range("K97").Value = 1404.1428
$ws.Range("M97").Value = -908.1428000000001
$ws.Range("H116").Value = 2087.2222
$ws.Range("I116").Value = 1338.0834
$ws.Range("K116").Value = 1338.0834
$ws.Range("M116").Value = 955.9166
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2087.2222
$ws.Range("I3").Value = 1338.0834
$ws.Range("K3").Value = 1338.0834
$ws.Range("M3").Value = -1224.0834
$ws.Range("H99").Value = 3588.7778
$ws.Range("I99").Value = 3537.375
$ws.Range("K99").Value = 3537.375
$ws.Range("M99").Value = -2039.375
$ws.Range("H134").Value = 14106.286
$ws.Range("I134").Value = 14106.286
$ws.Range("K134").Value = 42318.858
$ws.Range("M134").Value = -39783.858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6470.4287
$ws.Range("J31").Value = 7324.25
$ws.Range("L31").Value = 7324.25
$ws.Range("N31").Value = -7914.25
$ws.Range("H34").Value = 6470.4287
$ws.Range("J34").Value = 7324.25
$ws.Range("L34").Value = 7324.25
$ws.Range("N34").Value = -7728.25
$ws.Range("H99").Value = 5349.75
$ws.Range("I99").Value = 6466.6665
$ws.Range("K99").Value = 6466.6665
$ws.Range("M99").Value = -4968.6665
$ws.Range("H107").Value = 753
$ws.Range("I107").Value = 645.75
$ws.Range("J107").Value = 896
$ws.Range("K107").Value = 645.75
$ws.Range("L107").Value = 896
$ws.Range("M107").Value = 1274.25
$ws.Range("N107").Value = -4736
$ws.Range("H126").Value = 5349.75
$ws.Range("I126").Value = 6466.6665
$ws.Range("K126").Value = 19399.9995
$ws.Range("M126").Value = -16929.9995
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 1237
$ws.Range("I132").Value = 1237
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3711
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1181
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 9998.799999999999
$ws.Range("I134").Value = 6666.6665
$ws.Range("K134").Value = 19999.9995
$ws.Range("M134").Value = -17464.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 268.625
$ws.Range("I2").Value = 221.72728
$ws.Range("J2").Value = 371.8
$ws.Range("K2").Value = 1330.36368
$ws.Range("L2").Value = 2230.8
$ws.Range("M2").Value = -1217.36368
$ws.Range("N2").Value = -2456.8
$ws.Range("H131").Value = 13712.277
$ws.Range("I131").Value = 30285.285
$ws.Range("J131").Value = 3165.818
$ws.Range("K131").Value = 90855.855
$ws.Range("L131").Value = 9497.454000000002
$ws.Range("M131").Value = -85815.855
$ws.Range("N131").Value = -19577.454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14800
$ws.Range("H70").Value = 6702.4287
$ws.Range("I70").Value = 5966.3335
$ws.Range("J70").Value = 7254.5
$ws.Range("K70").Value = 5966.3335
$ws.Range("L70").Value = 7254.5
$ws.Range("M70").Value = -5696.3335
$ws.Range("N70").Value = -7794.5
$ws.Range("H73").Value = 6702.4287
$ws.Range("I73").Value = 5966.3335
$ws.Range("J73").Value = 7254.5
$ws.Range("K73").Value = 5966.3335
$ws.Range("L73").Value = 7254.5
$ws.Range("M73").Value = -5030.3335
$ws.Range("N73").Value = -9126.5
$ws.Range("H80").Value = 4198.75
$ws.Range("J80").Value = 4198.75
$ws.Range("L80").Value = 4198.75
$ws.Range("N80").Value = -6194.75
$ws.Range("H83").Value = 4198.75
$ws.Range("J83").Value = 4198.75
$ws.Range("L83").Value = 20993.75
$ws.Range("N83").Value = -30977.75
$ws.Range("H132").Value = 4498
$ws.Range("J132").Value = 4497.5
$ws.Range("L132").Value = 13492.5
$ws.Range("N132").Value = -18552.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3405.923
$ws.Range("I16").Value = 3514.75
$ws.Range("J16").Value = 2100
$ws.Range("K16").Value = 3514.75
$ws.Range("L16").Value = 2100
$ws.Range("M16").Value = -3344.75
$ws.Range("N16").Value = -2440
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("M27").ClearContents()
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10722
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12496
$ws.Range("M85").ClearContents()
$ws.Range("H100").Value = 6059.25
$ws.Range("I100").Value = 6414.5
$ws.Range("K100").Value = 6414.5
$ws.Range("M100").Value = -5873.5
$ws.Range("H122").Value = 4309.5713
$ws.Range("I122").Value = 4273.5
$ws.Range("K122").Value = 12820.5
$ws.Range("M122").Value = -10370.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9447
$ws.Range("I81").Value = 894.5
$ws.Range("J81").Value = 17999.5
$ws.Range("K81").Value = 1789
$ws.Range("L81").Value = 35999
$ws.Range("M81").Value = -728
$ws.Range("N81").Value = -38121
$ws.Range("H84").Value = 9447
$ws.Range("I84").Value = 894.5
$ws.Range("J84").Value = 17999.5
$ws.Range("K84").Value = 8945
$ws.Range("L84").Value = 179995
$ws.Range("M84").Value = -3641
$ws.Range("N84").Value = -190603
$ws.Range("H96").Value = 925
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("H107").Value = 360.75
$ws.Range("I107").Value = 347.66666
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1042.99998
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 877.0000199999999
$ws.Range("N107").Value = -5040
$ws.Range("H122").Value = 1574
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 2946
$ws.Range("I132").Value = 2966.8572
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 8900.571599999999
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -6370.571599999999
$ws.Range("N132").Value = -13460
